$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45171 -> 45172, i.e. 2023-09-02 -> 2023-09-03) for every data row (2-170).
$ws.Range("C2:C170").Value = 45172
